$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = -20.48559999999997
$ws.Range("A13").Value = -21.92210000000001
$ws.Range("A16").Value = -20.03319999999999
$ws.Range("A18").Value = -22.91530000000002
$ws.Range("A20").Value = -22.12950000000002
